# Updates the cryptos price list (columns D = Price, E = Volume(1h))
# with refreshed values, matching the "Updated cryptos list" GitHub Action
# commit. Values are written as text: for cells where the new text looks
# like a plain number we temporarily force a Text number format so Excel
# does not silently convert the cell to a numeric value, then clear the
# format again so the cell's style index is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '20.509.33'
$ws.Range('E2').Value = '  +2.15%  '

$ws.Range('D3').Value = '1.472.22'
$ws.Range('E3').Value = '  +3.56%  '

$ws.Range('E4').Value = '  +0.72%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9572'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -4.06%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '276.79'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.03%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3651'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.47%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3054'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.90%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '39.67'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.13%  '

$ws.Range('E10').Value = '  -0.64%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06605'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.79%  '

$ws.Range('E12').Value = '  +0.29%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.10'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.85%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.448'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.74%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.176'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.47%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001027'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.48%  '

$ws.Range('D17').Value = '1.474.75'
$ws.Range('E17').Value = '  +3.60%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.05892'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.09%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9636'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.45%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.04'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.83%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.456'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.84%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.48'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.62%  '

$ws.Range('E23').Value = '  -0.60%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.252'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.12%  '

$ws.Range('D25').Value = '20.561.98'
$ws.Range('E25').Value = '  +2.16%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '140.84'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +4.71%  '

$ws.Range('E27').Value = '  -7.64%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.16'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.99%  '

$ws.Range('D29').Value = '1.631.96'
$ws.Range('E29').Value = '  +3.02%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '113.50'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.03%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.952'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.50%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.8179'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.28%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.947'
$ws.Range('D33').ClearFormats()

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07927'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.73%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.539'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.82%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.233'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +11.09%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05764'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.62%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.722'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.97%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02034'
$ws.Range('D39').ClearFormats()

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '10.43'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.75%  '

$ws.Range('E41').Value = '  -4.26%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9563'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -4.09%  '

$ws.Range('E43').Value = '  -0.02%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5266'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.58%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.509'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.22%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.02'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.60%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '117.66'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.44%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5189'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.09%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.779'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.61%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06452'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.36%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9961'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.27%  '
